# Update odds data for the Japan J1 League matches (rows 3-5) and remove
# the USA - MLS match row (row 7), which has been dropped from this
# week's sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Kyoto vs Kawasaki Frontale) ---
$ws.Range("G3").Value = 2.55
$ws.Range("I3").Value = 2.5
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2.4
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 19
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 5.5
$ws.Range("Q3").Value = 1.53
$ws.Range("R3").Value = 2.4
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 3.75
$ws.Range("U3").Value = 1.44
$ws.Range("V3").Value = 2.63
$ws.Range("Z3").Value = 26
$ws.Range("AA3").Value = 17
$ws.Range("AC3").Value = 19
$ws.Range("AD3").Value = 7.5
$ws.Range("AF3").Value = 29
$ws.Range("AH3").Value = 13
$ws.Range("AT3").Value = 3.75
$ws.Range("AW3").Value = 5
$ws.Range("AY3").Value = 17

# --- Row 4 (Machida vs FC Tokyo) ---
$ws.Range("H4").Value = 3.6
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.93
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("W4").Value = 7.5
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 7
$ws.Range("AQ4").Value = 29
$ws.Range("AT4").Value = 3

# --- Row 5 (Sagan Tosu vs Yokohama F. Marinos) ---
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 4.1
$ws.Range("I5").Value = 1.95
$ws.Range("J5").Value = 3.5
$ws.Range("L5").Value = 2.4
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 26
$ws.Range("S5").Value = 1.2
$ws.Range("T5").Value = 4.33
$ws.Range("Y5").Value = 13
$ws.Range("AC5").Value = 26
$ws.Range("AO5").Value = 17
$ws.Range("AP5").Value = 19
$ws.Range("AQ5").Value = 51
$ws.Range("AT5").Value = 4.33

# --- Remove row 7 (Los Angeles FC vs Vancouver Whitecaps - USA MLS) ---
$ws.Rows.Item(7).Delete()
